$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.733.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.32%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.160.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.59%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.629'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.76%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.42'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.74%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.392'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0846'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.35%  '

# Row 11
$ws.Range("E11").Value = '  -0.14%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.93'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.97%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.479.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.53%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.86%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.804'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.35%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.89%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.162.20'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.25%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.642.32'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.02%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.59%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0845'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.38%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.96%  '

# Row 23
$ws.Range("E23").Value = '  +0.11%  '

# Row 24
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("E24").Value = '  +1.06%  '

# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.68%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.26%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.15%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.139'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.21%  '

# Row 29
$ws.Range("E29").Value = '  +3.53%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.68%  '

# Row 31
$ws.Range("E31").Value = '  +4.35%  '

# Row 32
$ws.Range("E32").Value = '  +1.36%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.59'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.34%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.67'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.72%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.95'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.88%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0619'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.27%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.66%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.60'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.64%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +20.82%  '

# Row 40
$ws.Range("E40").Value = '  -0.25%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '103.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.11%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0226'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.79%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.85%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.514.66'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.61%  '

# Row 45
$ws.Range("E45").Value = '  +0.96%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.06%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0920'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.16%  '

# Row 48
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.11%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.70%  '

# Row 50
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '50.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.65%  '

# Row 51
$ws.Range("E51").Value = '  +0.78%  '
